$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(111).Insert()

$ws.Cells.Item(111, 1).Value = 10
$ws.Cells.Item(111, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value = "La Araucanía"
$ws.Cells.Item(111, 4).Value = "2022-10-05"
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100107
$ws.Cells.Item(111, 8).Value = "Otros"
$ws.Cells.Item(111, 9).Value = 100107002
$ws.Cells.Item(111, 10).Value = "Chirimoya"
$ws.Cells.Item(111, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 25
$ws.Cells.Item(111, 14).Value = 3500
$ws.Cells.Item(111, 15).Value = 3500
$ws.Cells.Item(111, 16).Value = 3500
$ws.Cells.Item(111, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(111, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(111, 19).Value = 3500
$ws.Cells.Item(111, 20).Value = 1
